# Generate Report for Handoff
#
# The previous handoff ("ddbbd4ca-...md") finished and a new handoff
# ("5e667d6f-...md") plus a brand-new not-yet-localized file
# ("ffff67e95264-...md") now show up in the report, ahead of the
# always-present ".localization-config" row.

$wb = $excel.ActiveWorkbook

$oldFile    = "ddbbd4ca-3d83-49a6-9d26-3bb385f83087.md"
$readyFile  = "5e667d6f-dc90-457e-b295-fbe6abeb0028.md"
$newFile    = "ffff67e95264-01fd-41ee-a1f9-4aceb53d1d84.md"
$config     = ".localization-config"

$readyStatus  = "Ready for handoff"
$notLocalized = "Not to be localized"

$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/acdcb4c2f606574189073052b74560cd3051d4ad/.localization-config"
$readyUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/c3dd77abf733647f0bf1915376218c26f11ffc6a/e2e/$readyFile"
$newUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/c3dd77abf733647f0bf1915376218c26f11ffc6a/e2e/$newFile"

# ---- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Make room for the new "ffff67e95264..." row; the existing
# ".localization-config" row slides from row 3 down to row 4.
$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = $readyFile
$ws.Range("B2").Value = $readyStatus
$ws.Range("C2").Value = $readyStatus

$ws.Range("A3").Value = $newFile
$ws.Range("B3").Value = $readyStatus
$ws.Range("C3").Value = $readyStatus

# Row insertion does not carry hyperlinks along with it, so rebuild them.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $readyUrl,  "", "", $readyFile) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), $newUrl,    "", "", $newFile)   | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $config)    | Out-Null

# ---- per-locale sheets (zh-cn / de-de) -------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Handoff = "2016-01-28 05:51:49" },
    @{ Sheet = "de-de"; Handoff = "2016-01-28 05:51:59" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)
    $xlf = "5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec." + $loc.Sheet + ".xlf"
    $xlfUrl = "https://github.com/OpenLocalizationTest/oltest/blob/3821691a11d655d326ebf1527d572680262db9ec/targets/$xlf"

    $ws.Rows.Item(3).Insert()

    $ws.Range("A2").Value = $readyFile
    $ws.Range("B2").Value = $readyStatus
    $ws.Range("C2").Value = $xlf
    $ws.Range("D2").Value = $loc.Handoff
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Include"

    $ws.Range("A3").Value = $newFile
    $ws.Range("B3").Value = $readyStatus
    $ws.Range("C3").Value = $xlf
    $ws.Range("D3").Value = $loc.Handoff
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Include"

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $readyUrl,  "", "", $readyFile) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $xlfUrl,    "", "", $xlf)       | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $newUrl,    "", "", $newFile)   | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $xlfUrl,    "", "", $xlf)       | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $config)    | Out-Null
}

Write-Output "Generated handoff report"
